# Add a new "2020" column (Q) to the maternal-mortality-rate table, mirroring
# the existing year columns (D..P). For each data row we copy the formatting
# from the adjacent column P cell (so number format / font / border match)
# and then write the new value, and finally move the active selection like
# the source workbook did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (thin header separator row, no value) - just inherit P3's format.
$ws.Cells.Item(3, 16).Copy($ws.Cells.Item(3, 17))

# Row 4 - the new year header, 2020.
$ws.Cells.Item(4, 16).Copy($ws.Cells.Item(4, 17))
$ws.Cells.Item(4, 17).Value = 2020

# Data rows 5-14 - the new values for 2020.
$newValues = @{
    5  = 38.6
    6  = 42.4
    7  = 53.2
    8  = 90.6
    9  = 52.6
    10 = 24.5
    11 = 69.1
    12 = 32.2
    13 = 19.1
    14 = 25.2
}

foreach ($row in 5..14) {
    $ws.Cells.Item($row, 16).Copy($ws.Cells.Item($row, 17))
    $ws.Cells.Item($row, 17).Value = $newValues[$row]
}

# Move the active selection, matching the edited workbook.
$ws.Range("R27").Select()
